$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the overlap distance input value (I7): 29 -> 32
$ws.Range("I7").Value = 32

# Update the selected cell to I8
$ws.Range("I8").Select()
